$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.237798548211589, 2.219073147148707, 15.78507051212126, 3.973042978891778, 4.009743767593932, 46)
    3  = @(0.1668817152522232, 1.721208616355125, 11.78773906116087, 3.433327694986436, 3.468019574832006, 45)
    4  = @(0.007374364274256014, 1.802371724810875, 11.19060025837959, 3.345235456343782, 3.383901761293858, 44)
    5  = @(0.2454644955640485, 1.75127826504232, 10.92499624066271, 3.305298207524204, 3.335180323351736, 43)
    6  = @(0.08964530462448218, 1.817368728623215, 10.89571318232202, 3.300865520181339, 3.339645190262546, 42)
    7  = @(0.2922757529429554, 1.817158715501804, 11.47184296907565, 3.387010919538886, 3.416295991223567, 41)
    8  = @(0.1800782816237511, 1.881494406550722, 11.31509423920953, 3.363791646224471, 3.401759130395282, 40)
    9  = @(0.3440619902521568, 1.826509698037073, 11.52799779226949, 3.395290531349193, 3.421969063912408, 39)
    10 = @(0.1457618719346943, 1.855812420463658, 11.69936450046915, 3.420433379042655, 3.463198346522003, 38)
    11 = @(0.3616810767302359, 1.786912209159727, 11.95894640054174, 3.45817096172843, 3.486644870762215, 37)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
